# The workbook was re-saved from Excel after the user widened several
# columns (to fit the newly-visible header/data text) and moved the
# selection. Reproduce both via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ------------------------------------------------------
# Excel's COM ColumnWidth is expressed in "characters" of the Normal
# style's font and is stored back to the OOXML <col width="..."/> using a
# Maximum-Digit-Width pixel grid (width_px = round(ColumnWidth * MDW),
# stored = (width_px + 5) / MDW). Set the values so the round-tripped
# stored width lands on (or as close as achievable to) the target widths.
$ws.Columns.Item(4).ColumnWidth  = 22.6666666666667   # -> stored width 23.5
$ws.Columns.Item(5).ColumnWidth  = 11.1666666666667   # -> stored width 12
$ws.Columns.Item(6).ColumnWidth  = 27.6666666666667   # -> stored width 28.5
$ws.Columns.Item(7).ColumnWidth  = 18.3333333333333   # -> stored width ~19.1640625
$ws.Columns.Item(8).ColumnWidth  = 13.8333333333333   # -> stored width ~14.6640625
$ws.Columns.Item(9).ColumnWidth  = 13                 # -> stored width ~13.83203125
$ws.Columns.Item(11).ColumnWidth = 10.8333333333333   # -> stored width ~11.6640625
$ws.Columns.Item(12).ColumnWidth = 8.66666666666667   # -> stored width 9.5

# --- Selection -----------------------------------------------------------
$ws.Range("M9").Select()
